$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.24711565198759189
    2  = -0.0059999999422117867
    3  = -0.0039999999475881509
    4  = -0.0079999999045998749
    5  = -0.0029999999436389757
    6  = -0.0019999999362347864
    7  = -0.0099999998627082753
    8  = -0.0099999998627926523
    9  = -0.001999999937397412
    10 = -0.0019999999392403822
    11 = 0.047267379633091799
    12 = -0.0034999999239953006
    13 = -0.01121251193353956
    14 = -0.0079999998743218725
    15 = -0.00099999993738908444
    16 = -0.0019999999273072611
    17 = -0.0019999999259354695
    18 = -0.0039999999072559689
    19 = -0.0039999999589217516
    20 = -0.0039999999420121668
    21 = -0.055993854520801634
    22 = -0.0039999999380118112
    23 = -0.0049999999355083702
    24 = -0.019999999792929657
    25 = -0.019999999790274892
    26 = -0.0024999999289789798
    27 = -0.002499999925241525
    28 = -0.0019999999141315783
    29 = -0.0069999998555649867
    30 = -0.059999999364644285
    31 = -0.0069999998451528711
    32 = -0.0099999998170208215
    33 = -0.0039999998716755414
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
